# edit.ps1 - apply the "Store data in MySql and Accessing it is next" commit
#
# Summary of the change (per the OOXML diff):
#   1. Three paragraphs get a word ("analyzing" / "Analyze" / "BeautifulSoup")
#      wrapped in <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>
#      which also splits the surrounding text into extra <w:r> runs.
#   2. One paragraph ("... Clients access and view financial insights.")
#      additionally gets "Clients" wrapped in
#      <w:proofErr w:type="gramStart"/> ... <w:proofErr w:type="gramEnd"/>.
#   3. Many list-item paragraphs have an (originally separate) <w:r><w:tab/></w:r>
#      run merged together with the following text run into a single run.
#
# We implement every edit the same way: locate the target paragraph, then
# replace the run content of that paragraph (but not its paragraph mark) with
# freshly authored OOXML via Range.InsertXML - this lets us place <w:proofErr/>
# siblings and split/merge <w:r> runs exactly as the target document does.

$d = $word.ActiveDocument

function Find-ParagraphByText {
    param($doc, [string]$needle)
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    throw "Paragraph containing '$needle' not found"
}

function Set-ParagraphRunXml {
    # Replaces the run content of $para (everything up to, but not including,
    # the trailing paragraph mark) with the supplied inner OOXML fragment.
    param($doc, $para, [string]$innerXml)

    $startPos = $para.Range.Start
    $endPos = $para.Range.End - 1
    $target = $doc.Range($startPos, $endPos)

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
        + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
        + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
        + '<pkg:xmlData>' `
        + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
        + '<w:body><w:p>' + $innerXml + '</w:p></w:body>' `
        + '</w:document>' `
        + '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($pkg)
}

$rPr24 = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

# ---------------------------------------------------------------------------
# 1. Overview paragraph: split around "analyzing"
# ---------------------------------------------------------------------------
$pOverview = Find-ParagraphByText $d "project aims to provide an automated"
$innerOverview = `
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve">The &quot;Finance Insights Dashboard&quot; project aims to provide an automated, comprehensive solution for extracting, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr24 + '<w:t>analyzing</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr24 + '<w:t>, and presenting financial data in a user-friendly web interface. By leveraging web scraping techniques, data analysis, and interactive visualizations, this project delivers timely financial insights tailored to client needs.</w:t></w:r>'
Set-ParagraphRunXml $d $pOverview $innerOverview

# ---------------------------------------------------------------------------
# 2. "Data Analysis and Visualization" objective: split around "Analyze"
# ---------------------------------------------------------------------------
$pDataAnalysis = Find-ParagraphByText $d "Data Analysis and Visualization:"
$innerDataAnalysis = `
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve">- Data Analysis and Visualization: </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr24 + '<w:t>Analyze</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> the extracted data and create insightful visualizations to highlight key financial trends and metrics.</w:t></w:r>'
Set-ParagraphRunXml $d $pDataAnalysis $innerDataAnalysis

# ---------------------------------------------------------------------------
# 3. Plain tab+text run merges (no proofErr involved)
# ---------------------------------------------------------------------------
$tabMerges = @(
    @{ Needle = "Manage user settings";                         Lead = " ";     Text = "    - Manage user settings and preferences." },
    @{ Needle = "Define roles for admin and clients";           Lead = " ";     Text = "    - Define roles for admin and clients." },
    @{ Needle = "Admins manage the system and client access";   Lead = " ";     Text = "    - Admins manage the system and client access." },
    @{ Needle = "Develop scripts to extract financial data";    Lead = "    ";  Text = " - Develop scripts to extract financial data from one selected website." },
    @{ Needle = "Design a SQL database schema";                 Lead = "   ";   Text = "  - Design a SQL database schema to store the scraped data." },
    @{ Needle = "Implement data insertion, updating";           Lead = "   ";   Text = "  - Implement data insertion, updating, and retrieval mechanisms." },
    @{ Needle = "Clean and preprocess the scraped data";        Lead = "  ";    Text = "   - Clean and preprocess the scraped data using Pandas." },
    @{ Needle = "Handle missing values, inconsistent formats";  Lead = "  ";    Text = "   - Handle missing values, inconsistent formats, and duplicates." },
    @{ Needle = "Perform basic analysis to identify";           Lead = "  ";    Text = "   - Perform basic analysis to identify key financial metrics and trends." },
    @{ Needle = "Calculate metrics such as averages";           Lead = "   ";   Text = "  - Calculate metrics such as averages, growth rates, and comparisons." },
    @{ Needle = "Create key interactive charts";                Lead = "   ";   Text = "  - Create key interactive charts using D3.js or Chart.js." },
    @{ Needle = "Types of visualizations may include";          Lead = "   ";   Text = "  - Types of visualizations may include line charts, bar charts, and pie charts." },
    @{ Needle = "Generate summary reports displaying";          Lead = "   ";   Text = "  - Generate summary reports displaying key financial insights." },
    @{ Needle = "Implement the backend using Flask";            Lead = "   ";   Text = "  - Implement the backend using Flask to serve data and handle user requests." },
    @{ Needle = "Develop API endpoints for data retrieval";     Lead = "   ";   Text = "  - Develop API endpoints for data retrieval and user management." },
    @{ Needle = "Create a responsive web interface using HTML"; Lead = "     "; Text = "- Create a responsive web interface using HTML, CSS, and JavaScript."; NoTrailingSpace = $true },
    @{ Needle = "Design intuitive navigation and layouts";      Lead = "    "; Text = " - Design intuitive navigation and layouts for easy data interpretation." },
    @{ Needle = "Embed interactive visualizations into";        Lead = "    "; Text = " - Embed interactive visualizations into the web interface." },
    @{ Needle = "Ensure seamless interaction between";          Lead = "   ";  Text = "  - Ensure seamless interaction between frontend and backend components." },
    @{ Needle = "Notify users about the latest data updates";   Lead = "   ";  Text = "  - Notify users about the latest data updates and reports." },
    @{ Needle = "Provide alerts for system issues";             Lead = " ";    Text = "    - Provide alerts for system issues or updates to users and admins." }
)

foreach ($item in $tabMerges) {
    $p = Find-ParagraphByText $d $item.Needle
    $leadXml = '<w:r><w:t xml:space="preserve">' + $item.Lead + '</w:t></w:r>'
    if ($item.NoTrailingSpace) {
        $textXml = '<w:r><w:tab/><w:t>' + $item.Text + '</w:t></w:r>'
    } else {
        $textXml = '<w:r><w:tab/><w:t xml:space="preserve">' + $item.Text + '</w:t></w:r>'
    }
    Set-ParagraphRunXml $d $p ($leadXml + $textXml)
}

# ---------------------------------------------------------------------------
# 4. "Clients access and view financial insights." - tab merge AND gramStart/End
# ---------------------------------------------------------------------------
$pClients = Find-ParagraphByText $d "Clients access and view financial insights"
$innerClients = `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:tab/><w:t xml:space="preserve">    - </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>Clients</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> access and view financial insights.</w:t></w:r>'
Set-ParagraphRunXml $d $pClients $innerClients

# ---------------------------------------------------------------------------
# 5. Technology Used / Python paragraph: split around "BeautifulSoup"
# ---------------------------------------------------------------------------
$pPython = Find-ParagraphByText $d "Employed for web scraping and data analysis"
$innerPython = `
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve">- Python: Employed for web scraping and data analysis. Libraries such as </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr24 + '<w:t>BeautifulSoup</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> and Scrapy will be used for scraping, while Pandas and Matplotlib/Seaborn will handle data analysis and visualization.</w:t></w:r>'
Set-ParagraphRunXml $d $pPython $innerPython

Write-Output "Applied all edits"
